$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 33; $row++) {
    $cell = $ws.Cells.Item($row, 15)  # column O is the 15th column
    if ($cell.Value2 -eq "2022-08-02 14:44:24") {
        $cell.Value2 = "2022-08-02 20:56:57"
    }
}
